$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Daily roll-forward over the data rows (row 1 is the header):
#   - "剩余" (E, days remaining) counts down by 1 each day;
#   - once it would drop below 1 it resets to 10 and "开始时间"
#     (F, a yyyyMMdd-style serial) jumps forward 10 days to start a
#     new 10-day cycle.
# Rows whose start-date cell isn't a well-formed 8-digit 2026 date
# (a bad data entry) are left untouched, same as upstream.
for ($row = 2; $row -le 99; $row++) {
    $eCell = $ws.Cells.Item($row, 5)
    $fCell = $ws.Cells.Item($row, 6)

    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    if ($eVal -eq $null -or $fVal -eq $null) { continue }

    $fText = [string]$fVal
    if ($fText.Length -ne 8 -or -not $fText.StartsWith("2026")) { continue }

    if ($eVal -eq 1) {
        $eCell.Value2 = 10
        $fCell.Value2 = $fVal + 10
    } else {
        $eCell.Value2 = $eVal - 1
    }
}
